$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Formula = "=B2/60"

$ws.Range("C3").Select() | Out-Null
